# RDBES Data Model workbook edit
# Commit message: "SelectionMethod is replacing RS_SelectionMethod. V19.1.2"
#
# This script:
#  1) Replaces every occurrence of the vocabulary code-list name
#     "RS_SelectionMethod" with its new name "SelectionMethod" across the
#     sheets that reference it in their "Vocabulary code type" (K) column.
#  2) Renames the main model sheet from "Model v1.19.1" to "Model v1.19.2".
#  3) Moves the active sheet/tab selection from "Onshore Event" to the main
#     "Model v1.19.2" sheet.
#  4) Restores the cursor/selection position on a couple of sheets that were
#     left scrolled/selected elsewhere by the author.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename "SelectionMethod" code list references.
# ---------------------------------------------------------------------------
$cellsToRename = @(
    @{ Sheet = "Temporal Event";        Cells = @("K18","K20") },
    @{ Sheet = "Location";              Cells = @("K18","K20") },
    @{ Sheet = "Vessel Selection";      Cells = @("K18","K20") },
    @{ Sheet = "Fishing Trip";          Cells = @("K29","K31") },
    @{ Sheet = "Fishing Operation";     Cells = @("K51","K53") },
    @{ Sheet = "Onshore Event";         Cells = @("K22","K24") },
    @{ Sheet = "Landing event";         Cells = @("K46","K48") },
    @{ Sheet = "Species Selection";     Cells = @("K25","K27") },
    @{ Sheet = "Sample";                Cells = @("K38") },
    @{ Sheet = "Biological Variable";   Cells = @("K25") }
)

foreach ($entry in $cellsToRename) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    foreach ($addr in $entry.Cells) {
        $rng = $ws.Range($addr)
        if ($rng.Value2 -eq "RS_SelectionMethod") {
            $rng.Value2 = "SelectionMethod"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Rename the main sheet to the new version number.
# ---------------------------------------------------------------------------
$modelSheet = $wb.Worksheets.Item("Model v1.19.1")
$modelSheet.Name = "Model v1.19.2"

# ---------------------------------------------------------------------------
# 3) Re-point the active tab at the renamed main sheet.
# ---------------------------------------------------------------------------
$modelSheet.Activate()
$modelSheet.Range("A7").Select()

# ---------------------------------------------------------------------------
# 4) Restore the selections left on a couple of sheets.
# ---------------------------------------------------------------------------
$temporalEvent = $wb.Worksheets.Item("Temporal Event")
$temporalEvent.Activate()
$temporalEvent.Application.ActiveWindow.ScrollRow = 10
$temporalEvent.Range("E15").Select()

$frequencyMeasure = $wb.Worksheets.Item("Frequency Measure")
$frequencyMeasure.Activate()
$frequencyMeasure.Range("K9").Select()

$biologicalVariable = $wb.Worksheets.Item("Biological Variable")
$biologicalVariable.Activate()
$biologicalVariable.Range("K13").Select()

# Finally re-activate the main model sheet so it is the one shown/selected
# when the workbook is saved.
$modelSheet.Activate()
